$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.459.83'
$ws.Range("E2").Value = '  -0.80%  '
$ws.Range("D3").Value = '1.566.36'
$ws.Range("E3").Value = '  -1.08%  '
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '208.46'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.07%  '
$ws.Range("E7").Value = '  -0.29%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.09'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.64%  '
$ws.Range("E9").Value = '  -1.41%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0592'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.09%  '
$ws.Range("E11").Value = '  -0.29%  '
$ws.Range("D12").Value = '1.787.71'
$ws.Range("E12").Value = '  -1.10%  '
$ws.Range("D13").Value = '1.588.53'
$ws.Range("E13").Value = '  +0.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.82'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.519'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.66'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.70%  '
$ws.Range("D17").Value = '27.446.33'
$ws.Range("E17").Value = '  -0.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '213.07'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.04%  '
$ws.Range("E19").Value = '  -0.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.83%  '
$ws.Range("E21").Value = '  -0.27%  '
$ws.Range("E22").Value = '  -0.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.56'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.75%  '
$ws.Range("E24").Value = '  +3.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.37'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("E26").Value = '  -0.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.68'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.10%  '
$ws.Range("E28").Value = '  -1.00%  '
$ws.Range("E29").Value = '  -1.71%  '
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("E31").Value = '  +1.14%  '
$ws.Range("E32").Value = '  -0.56%  '
$ws.Range("D33").Value = '1.375.57'
$ws.Range("E33").Value = '  -0.77%  '
$ws.Range("E34").Value = '  +1.30%  '
$ws.Range("E35").Value = '  +0.89%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.958'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.31%  '
$ws.Range("E37").Value = '  -0.99%  '
$ws.Range("E38").Value = '  +1.04%  '
$ws.Range("E39").Value = '  -1.80%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.823'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.975'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.16%  '
$ws.Range("E43").Value = '  +2.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '64.06'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.96%  '
$ws.Range("E45").Value = '  -0.30%  '
$ws.Range("E46").Value = '  +0.74%  '
$ws.Range("D47").Value = '1.701.56'
$ws.Range("E47").Value = '  -1.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.47'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.83%  '
$ws.Range("D49").Value = '0.0₆0100'
$ws.Range("E49").Value = '  -0.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0959'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.54%  '
$ws.Range("E51").Value = '  -0.63%  '
